# Auto-generated edit script: updates cryptos.xlsx Sheet1 price/volume data
# to reflect the scraped values from the Sat Jun  1 21:44:13 UTC 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.819.49'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '3.816.92'
$ws.Range('E3').Value = '  +0.90%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''606.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.95%  '
$ws.Range('D6').Value = '''167.22'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +0.47%  '
$ws.Range('E9').Value = '  +0.69%  '
$ws.Range('D10').Value = '''6.31'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.70%  '
$ws.Range('D11').Value = '''0.452'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.70%  '
$ws.Range('E12').Value = '  -0.85%  '
$ws.Range('D13').Value = '''36.08'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Value = '4.454.32'
$ws.Range('E14').Value = '  +0.81%  '
$ws.Range('D15').Value = '3.829.08'
$ws.Range('E15').Value = '  +0.76%  '
$ws.Range('D16').Value = '''18.53'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.56%  '
$ws.Range('D17').Value = '67.822.73'
$ws.Range('E17').Value = '  +0.39%  '
$ws.Range('E18').Value = '  +1.40%  '
$ws.Range('E19').Value = '  +0.40%  '
$ws.Range('D20').Value = '''463.08'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.48%  '
$ws.Range('D21').Value = '''9.94'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.37%  '
$ws.Range('E22').Value = '  +0.89%  '
$ws.Range('E23').Value = '  -2.67%  '
$ws.Range('D24').Value = '''83.44'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('D25').Value = '''12.09'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.55%  '
$ws.Range('D26').Value = '''2.12'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.14%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').Value = '''10.07'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D28').Value = '''1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('D29').Value = '3.966.23'
$ws.Range('E29').Value = '  +0.86%  '
$ws.Range('D30').Value = '''2.81'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.62%  '
$ws.Range('E31').Value = '  +1.88%  '
$ws.Range('D32').Value = '''2.23'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.64%  '
$ws.Range('D33').Value = '''29.67'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.49%  '
$ws.Range('E34').Value = '  +0.14%  '
$ws.Range('D35').Value = '''9.10'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.53%  '
$ws.Range('D36').Value = '3.757.33'
$ws.Range('E36').Value = '  +0.57%  '
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('D38').Value = '''3.38'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.57%  '
$ws.Range('E39').Value = '  +0.16%  '
$ws.Range('D40').Value = '''0.999'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.22%  '
$ws.Range('E41').Value = '  +0.72%  '
$ws.Range('D42').Value = '''0.999'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D44').Value = '''48.15'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.27%  '
$ws.Range('E45').Value = '  +0.71%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '''28.47'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +10.96%  '
$ws.Range('B47').Value = 'Arweave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D47').Value = '''43.19'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.82%  '
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').Value = '''1.38'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +10.90%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').Value = '''8.36'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.22%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').Value = '''148.79'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('E51').Value = '  +0.40%  '
